# Added one more RegistrationTest
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RegistrationUser
$ws2 = $wb.Worksheets.Item(2)   # LoginUser

# --- New test row on the RegistrationUser sheet -----------------------
# Columns: A=TestName, B=Email, C=FullName, D=Password, E=ConfirmPassword
# Write in this order so new shared strings land in the same slots as
# the authored workbook (16=RegistrationWithoutPassword, 17=IvanIvanov,
# 18=email@abv.bg).
$ws1.Range("A4").Value = "RegistrationWithoutPassword"
$ws1.Range("C4").Value = "IvanIvanov"
$ws1.Range("B4").Value = "email@abv.bg"

# --- Column width tweaks on RegistrationUser (no longer auto bestFit) --
$ws1.Columns.Item(2).ColumnWidth = 14.33203125
$ws1.Columns.Item(3).ColumnWidth = 12.5546875
$ws1.Columns.Item(4).ColumnWidth = 15.6640625

# --- Move the active sheet / selections --------------------------------
# LoginUser stops being the active tab, its lingering selection moves to C6
[void]$ws2.Activate()
[void]$ws2.Range("C6").Select()

# RegistrationUser becomes the active tab, selection moves to E6
[void]$ws1.Activate()
[void]$ws1.Range("E6").Select()
